$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Rename the shortened test case names in column A (rows 3-5)
$ws.Range("A3").Value = "02_Normal"
$ws.Range("A4").Value = "03_EstabEntry"
$ws.Range("A5").Value = "04_Tourism"

# Update the selection on the active sheet to A6
$ws.Range("A6").Select()
